# Add a new student row (row 5) to the "Danh Sách Nhóm 2 Web" roster, copying
# the row-3/4 formatting (bordered, centered, wrap-capable style) down onto it,
# and leave a stray value in E11 (matching the source commit), updating the
# sheet's used-range dimension and active selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing data row (row 4) onto the new
# row 5 before writing values into it, so the new row inherits the same
# cell style (borders, centered/wrapped alignment) as the rows above it.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

# New student entry
$ws.Range("A5").Value = 3122411119
$ws.Range("B5").Value = "Đặng Nguyễn Thành "
$ws.Range("C5").Value = "Luân"
$ws.Range("D5").Value = "DCT122C2"
$ws.Range("E5").Value = 2

# Match the wrapped-text row height used by the other data rows.
$ws.Rows.Item(5).RowHeight = 26.4

# Stray cell value far below the table (as in the source workbook).
$ws.Range("E11").Value = "s"

# Move/keep the active selection on E11, matching the saved view state.
$ws.Range("E11").Select()
